$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "29.900.92"
Set-TextCell 2 5 "  +1.16%  "
Set-TextCell 3 4 "1.628.81"
Set-TextCell 3 5 "  +1.85%  "
Set-TextCell 4 5 "  -0.02%  "
Set-TextCell 5 4 "214.42"
Set-TextCell 5 5 "  +1.13%  "
Set-TextCell 6 4 "0.521"
Set-TextCell 6 5 "  +1.47%  "
Set-TextCell 7 5 "  +0.02%  "
Set-TextCell 8 4 "29.39"
Set-TextCell 8 5 "  +9.61%  "
Set-TextCell 9 4 "0.259"
Set-TextCell 9 5 "  +3.31%  "
Set-TextCell 10 4 "0.0613"
Set-TextCell 10 5 "  +2.40%  "
Set-TextCell 11 4 "0.0916"
Set-TextCell 11 5 "  +0.78%  "
Set-TextCell 12 4 "1.863.21"
Set-TextCell 12 5 "  +2.01%  "
Set-TextCell 13 4 "1.642.59"
Set-TextCell 13 5 "  +2.77%  "
Set-TextCell 14 4 "0.569"
Set-TextCell 14 5 "  +6.21%  "
Set-TextCell 15 4 "3.90"
Set-TextCell 15 5 "  +4.72%  "
Set-TextCell 16 4 "29.932.07"
Set-TextCell 16 5 "  +1.26%  "
Set-TextCell 17 4 "9.00"
Set-TextCell 17 5 "  +18.98%  "
Set-TextCell 18 4 "64.81"
Set-TextCell 18 5 "  +1.86%  "
Set-TextCell 19 4 "246.19"
Set-TextCell 19 5 "  +2.70%  "
Set-TextCell 20 4 "0.0₃0703"
Set-TextCell 20 5 "  +1.44%  "
Set-TextCell 21 5 "  -0.02%  "
Set-TextCell 22 4 "4.13"
Set-TextCell 22 5 "  +3.61%  "
Set-TextCell 23 4 "9.57"
Set-TextCell 23 5 "  +3.55%  "
Set-TextCell 24 4 "2.11"
Set-TextCell 24 5 "  +1.14%  "
Set-TextCell 25 4 "158.46"
Set-TextCell 25 5 "  +2.36%  "
Set-TextCell 26 4 "15.68"
Set-TextCell 26 5 "  +2.18%  "
Set-TextCell 27 5 "  +2.16%  "
Set-TextCell 28 4 "6.60"
Set-TextCell 28 5 "  +3.17%  "
Set-TextCell 29 5 "  +0.01%  "
Set-TextCell 30 4 "0.0491"
Set-TextCell 30 5 "  +3.06%  "
Set-TextCell 31 4 "1.13"
Set-TextCell 31 5 "  +6.47%  "
Set-TextCell 32 5 "  +3.85%  "
Set-TextCell 33 4 "3.20"
Set-TextCell 33 5 "  +1.95%  "
Set-TextCell 34 4 "1.428.80"
Set-TextCell 34 5 "  -0.43%  "
Set-TextCell 35 4 "1.65"
Set-TextCell 35 5 "  +6.32%  "
Set-TextCell 36 5 "  +0.51%  "
Set-TextCell 37 5 "  +1.80%  "
Set-TextCell 38 5 "  -0.95%  "
Set-TextCell 39 5 "  +3.29%  "
Set-TextCell 40 4 "0.557"
Set-TextCell 40 5 "  +3.35%  "
Set-TextCell 41 2 "Kaspa"
Set-TextCell 41 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell 41 4 "0.0501"
Set-TextCell 41 5 "  +1.42%  "
Set-TextCell 42 2 "ARBITRUM"
Set-TextCell 42 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 42 4 "0.832"
Set-TextCell 42 5 "  +3.96%  "
Set-TextCell 43 5 "  +6.90%  "
Set-TextCell 44 2 "BitcoinSV"
Set-TextCell 44 3 "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextCell 44 4 "54.67"
Set-TextCell 44 5 "  +2.57%  "
Set-TextCell 45 2 "RenderToken"
Set-TextCell 45 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 45 4 "1.96"
Set-TextCell 45 5 "  +0.18%  "
Set-TextCell 46 4 "69.10"
Set-TextCell 46 5 "  +5.17%  "
Set-TextCell 47 5 "  +0.00%  "
Set-TextCell 48 4 "5.45"
Set-TextCell 48 5 "  +2.27%  "
Set-TextCell 49 4 "1.770.82"
Set-TextCell 49 5 "  +2.04%  "
Set-TextCell 50 4 "89.64"
Set-TextCell 50 5 "  +3.84%  "
Set-TextCell 51 5 "  +4.05%  "
